# Applies the "Data.xlsx" update:
#  - Trees sheet: correct Dingle's lat/long (row 2) to the more precise
#    fix, and log a new "example" tree with a "notes" start-of-season note.
#  - Saps sheet: log two new sap entries (Dingle + example tree).
#  - Syrups sheet: log a new syrup batch.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Trees
# ---------------------------------------------------------------------
$trees = $wb.Worksheets.Item("Trees")

# Row 2 (Dingle) gets a corrected GPS fix.
$trees.Range("F2").Value = 39.4921
$trees.Range("G2").Value = -74.5323

# New row 7: an "example" tree entry.
$trees.Range("A7").Value = "example"
$trees.Range("B7").Value = 10
$trees.Range("C7").Value = 0
$trees.Range("D2").Copy($trees.Range("D7"))
$trees.Range("D7").Value = 44295
$trees.Range("E7").Value = 10
$trees.Range("F7").Value = 39.399
$trees.Range("G7").Value = -74.5146
$trees.Range("H7").Value = "notes"

# ---------------------------------------------------------------------
# Saps
# ---------------------------------------------------------------------
$saps = $wb.Worksheets.Item("Saps")

# New row 19: another Dingle sap entry.
$saps.Range("A19").Value = "Dingle"
$saps.Range("B19").Value = 55
$saps.Range("C2").Copy($saps.Range("C19"))
$saps.Range("C19").Value = 44295
$saps.Range("D19").Value = 56

# New row 20: the example tree's sap entry.
$saps.Range("A20").Value = "example"
$saps.Range("B20").Value = 10
$saps.Range("C2").Copy($saps.Range("C20"))
$saps.Range("C20").Value = 44295
$saps.Range("D20").Value = 56

# ---------------------------------------------------------------------
# Syrups
# ---------------------------------------------------------------------
$syrups = $wb.Worksheets.Item("Syrups")

# New row 5: a new syrup batch made with propane.
$syrups.Range("A5").Value = 3
$syrups.Range("B5").Value = 10
$syrups.Range("C5").Value = 5
$syrups.Range("D2").Copy($syrups.Range("D5"))
$syrups.Range("D5").Value = 44295
$syrups.Range("E5").Value = 2
$syrups.Range("F5").Value = 30
$syrups.Range("G5").Value = "Propane"
